$d = $word.ActiveDocument

# Locate the paragraph that ends with "... để tuân thủ Material 3."
$findRng = $d.Content
$found = $findRng.Find.Execute("để tuân thủ Material 3.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'để tuân thủ Material 3.'"
}

$anchorPara = $findRng.Paragraphs(1)

# Insert a brand-new (empty) paragraph right after the anchor paragraph.
# It inherits the anchor's paragraph formatting (numbering, spacing, rPr),
# which matches the bullet list (numId 45) the new item should belong to.
$anchorRng = $anchorPara.Range
$anchorRng.Collapse(0)
$anchorRng.InsertParagraphAfter()

# Grab the freshly created paragraph (immediately follows the anchor).
$newPara = $anchorPara.Next()
$newRng = $newPara.Range

# New list item text (a single run with no explicit run formatting,
# matching the target markup exactly).
$newText = "Don't use 'BuildContext's across async gaps, guarded by an unrelated 'mounted' check. Guard a 'State.context' use with a 'mounted' check on the State, and other BuildContext use with a 'mounted' check on the BuildContext."

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="45"/></w:numPr>
<w:spacing w:before="100" w:beforeAutospacing="1" w:after="0" w:line="240" w:lineRule="auto"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:color w:val="1B1C1D"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
</w:rPr>
</w:pPr>
<w:r>
<w:t>$newText</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $newRng.InsertXML($xml)
